# Update the NR yearly financials sheet: a new reporting period (column)
# is inserted immediately before column D, shifting the existing D:K data
# one column to the right (E:L), and the new column D is populated with
# the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NR")

# 1) Insert a new blank column before D; this shifts D:K -> E:L (values,
#    formulas and styles all move together).
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D has no number formatting / font of its
#    own yet. Column E now holds what used to live in D, so copy its
#    formatting back onto D (the two columns should look identical) and
#    seed D's values from E as a baseline - most rows repeat the same
#    figure (0 / "NA" / blank) across every period, so this gives every
#    row in D the correct starting value before the real updates below.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = 0

# The blank separator rows (36/78 - fully empty - and 37/79, which only
# carry a section-title label in column B) never had a column D cell
# before the insert; the bulk paste above touched every row in the
# 7:102 span, so undo that incidental spillover here.
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# 3) Write the new period's actual figures into column D.
$newValues = @{
    7   = 43465
    8   = 946500
    9   = 767000
    10  = 179600
    17  = 883000
    18  = 63600
    20  = -1400
    21  = 108000
    22  = 14900
    23  = 47300
    24  = 16600
    26  = 30700
    27  = 30700
    29  = 1600
    32  = 1400
    33  = 32300
    35  = 32300
    38  = 43465
    41  = 56100
    43  = 254400
    44  = 196900
    45  = 15900
    46  = 523300
    48  = 316300
    49  = 69000
    52  = 7300
    54  = 915900
    57  = 90600
    58  = 2500
    59  = 48800
    60  = 141900
    61  = 159200
    62  = 45000
    66  = 346200
    72  = 148800
    76  = 569700
    80  = 43465
    81  = 32300
    83  = 45900
    89  = 63400
    91  = -45100
    94  = -55800
    100 = -4500
    101 = -4300
    102 = -1200
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D" + $row).Value = $newValues[$row]
}
